# Update "想去人数" (number of people wanting to go) figures that changed
# between data snapshots, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F9").Value = 6265
$wsExhibit.Range("F21").Value = 4604
$wsExhibit.Range("F25").Value = 195

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 6265
$wsAll.Range("F21").Value = 4604
$wsAll.Range("F26").Value = 195
